$wb = $excel.ActiveWorkbook

# Add the new "ProductList" worksheet. Add() inserts it right after the
# currently active sheet ("Login"/"RegistrationDetails" active tab), which
# lands it between "Login" and "RegistrationDetails" - matching the target
# sheet order/ids.
$ws = $wb.Worksheets.Add()
$ws.Name = "ProductList"

# --- Cell values --------------------------------------------------------
# Written in the same order the shared strings were introduced in the
# target workbook: PRODUCTS, Forbidden Planet, Asus Eee PC Touch,
# EKTORP Neckroll, Quantity.
$ws.Range("B2").Value = "PRODUCTS"
$ws.Range("B3").Value = "Forbidden Planet"
$ws.Range("B4").Value = "Asus Eee PC Touch"
$ws.Range("B5").Value = "EKTORP Neckroll"
$ws.Range("C2").Value = "Quantity"
$ws.Range("C3").Value = 1
$ws.Range("C4").Value = 1
$ws.Range("C5").Value = 1

# --- Formatting ----------------------------------------------------------
# Build the bold / bordered / centered header look on a scratch cell first
# and paste the finished format in one shot, so only a single combined
# style entry is produced for the header row (instead of one new style per
# individual property write).
$scratch = $ws.Range("Z100")
$scratch.Font.Bold = $true
$scratch.HorizontalAlignment = -4131
$scratch.VerticalAlignment = -4108
$scratch.Borders.LineStyle = 1
$scratch.Copy()
$hdr = $ws.Range("B2:C2")
$hdr.PasteSpecial(-4122)
$ws.Range("Z100").Clear()

# Plain bordered look for the data rows.
$body = $ws.Range("B3:C5")
$body.Borders.LineStyle = 1

# --- Column widths ---------------------------------------------------
$ws.Columns("B").ColumnWidth = 17.5703125
$ws.Columns("C").ColumnWidth = 10.85546875

# --- Selection / active cell -----------------------------------------
$ws.Range("I12").Select()
